$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Buying Opportunity" column (B)
$ws.Range("B2").Value  = "NSE:CCL"
$ws.Range("B3").Value  = "NSE:CUPID"
$ws.Range("B4").Value  = "NSE:DEEPINDS"
$ws.Range("B5").Value  = "NSE:HINDWAREAP"
$ws.Range("B6").Value  = "NSE:KAPSTON"
$ws.Range("B7").Value  = "NSE:MAFANG"
$ws.Range("B8").Value  = "NSE:NIPPOBATRY"
$ws.Range("B9").Value  = "NSE:ORCHPHARMA"
$ws.Range("B10").Value = "NSE:PDMJEPAPER"

# Update the "support Zone" column (C), including new rows 16-28
$ws.Range("C2").Value  = "NSE:AARTIIND"
$ws.Range("C3").Value  = "NSE:BARBEQUE"
$ws.Range("C4").Value  = "NSE:BASF"
$ws.Range("C5").Value  = "NSE:CESC"
$ws.Range("C6").Value  = "NSE:COFFEEDAY"
$ws.Range("C7").Value  = "NSE:DAMODARIND"
$ws.Range("C8").Value  = "NSE:DENORA"
$ws.Range("C9").Value  = "NSE:DIAMONDYD"
$ws.Range("C10").Value = "NSE:FOSECOIND"
$ws.Range("C11").Value = "NSE:GARFIBRES"
$ws.Range("C12").Value = "NSE:GHCL"
$ws.Range("C13").Value = "NSE:GHCLTEXTIL"
$ws.Range("C14").Value = "NSE:GREENPLY"
$ws.Range("C15").Value = "NSE:GUFICBIO"
$ws.Range("C16").Value = "NSE:HAPPSTMNDS"
$ws.Range("C17").Value = "NSE:IMFA"
$ws.Range("C18").Value = "NSE:JHS"
$ws.Range("C19").Value = "NSE:JPOLYINVST"
$ws.Range("C20").Value = "NSE:MANAKSIA"
$ws.Range("C21").Value = "NSE:MAYURUNIQ"
$ws.Range("C22").Value = "NSE:MEDICAMEQ"
$ws.Range("C23").Value = "NSE:MOLDTKPAC"
$ws.Range("C24").Value = "NSE:MUKANDLTD"
$ws.Range("C25").Value = "NSE:NIACL"
$ws.Range("C26").Value = "NSE:PRIMESECU"
$ws.Range("C27").Value = "NSE:RELAXO"
$ws.Range("C28").Value = "NSE:RPGLIFE"

# "long buildup" column (D) - only row 2 changes
$ws.Range("D2").Value = "NSE:MARICO"

# "FII ENTERING" column (F) - row 2 changes, row 3 cleared
$ws.Range("F2").Value = "NSE:MARICO"
$ws.Range("F3").Value = ""

# "Short buildup" column (E) - rows 2-5 cleared out
$ws.Range("E2").Value = ""
$ws.Range("E3").Value = ""
$ws.Range("E4").Value = ""
$ws.Range("E5").Value = ""

# New index rows (column A) for the newly appended rows 16-28.
# Copy formatting (bold, border, centered) from the last existing index
# cell (A15) down through the new rows, then set the sequential values.
$ws.Range("A15").Copy($ws.Range("A16:A28"))

$ws.Range("A16").Value = 14
$ws.Range("A17").Value = 15
$ws.Range("A18").Value = 16
$ws.Range("A19").Value = 17
$ws.Range("A20").Value = 18
$ws.Range("A21").Value = 19
$ws.Range("A22").Value = 20
$ws.Range("A23").Value = 21
$ws.Range("A24").Value = 22
$ws.Range("A25").Value = 23
$ws.Range("A26").Value = 24
$ws.Range("A27").Value = 25
$ws.Range("A28").Value = 26
